$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# row 5
$ws.Range("H5").Value = 905.2857
$ws.Range("I5").Value = 259
$ws.Range("K5").Value = 259
$ws.Range("M5").Value = -144
# row 76
$ws.Range("H76").Value = 7777.1113
$ws.Range("J76").Value = 7999.2856
$ws.Range("L76").Value = 7999.2856
$ws.Range("N76").Value = -8629.285599999999
# row 79
$ws.Range("H79").Value = 7777.1113
$ws.Range("J79").Value = 7999.2856
$ws.Range("L79").Value = 7999.2856
$ws.Range("N79").Value = -10183.2856
# row 100
$ws.Range("H100").Value = 2148.7273
$ws.Range("I100").Value = 1634.1428
$ws.Range("J100").Value = 3049.25
$ws.Range("K100").Value = 1634.1428
$ws.Range("L100").Value = 3049.25
$ws.Range("M100").Value = -1093.1428
$ws.Range("N100").Value = -4131.25

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 5319.7
$ws.Range("I2").Value = 4537.8
$ws.Range("K2").Value = 4537.8
$ws.Range("M2").Value = -4424.8
# row 14
$ws.Range("H14").Value = 323.66666
$ws.Range("I14").Value = 188.4
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 188.4
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = -13.40000000000001
$ws.Range("N14").Value = -1350
# row 116
$ws.Range("H116").Value = 5319.7
$ws.Range("I116").Value = 4537.8
$ws.Range("K116").Value = 4537.8
$ws.Range("M116").Value = -2243.8
# row 122
$ws.Range("H122").Value = 20836442
$ws.Range("I122").Value = 2334
$ws.Range("K122").Value = 7002
$ws.Range("M122").Value = -4552
# row 132
$ws.Range("H132").Value = 18546510
$ws.Range("I132").Value = 2036.8276
$ws.Range("K132").Value = 6110.4828
$ws.Range("M132").Value = -3580.4828

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 5319.7
$ws.Range("I3").Value = 4537.8
$ws.Range("K3").Value = 4537.8
$ws.Range("M3").Value = -4423.8
# row 86
$ws.Range("H86").Value = 10665.125
$ws.Range("I86").Value = 14333.333
$ws.Range("J86").Value = 5948.857
$ws.Range("K86").Value = 14333.333
$ws.Range("L86").Value = 5948.857
$ws.Range("M86").Value = -13210.333
$ws.Range("N86").Value = -8194.857
# row 89
$ws.Range("H89").Value = 10665.125
$ws.Range("I89").Value = 14333.333
$ws.Range("J89").Value = 5948.857
$ws.Range("K89").Value = 71666.66500000001
$ws.Range("L89").Value = 29744.285
$ws.Range("M89").Value = -66050.66500000001
$ws.Range("N89").Value = -40976.285
# row 134
$ws.Range("H134").Value = 2419.282
$ws.Range("I134").Value = 2694.1333
$ws.Range("J134").Value = 1503.1111
$ws.Range("K134").Value = 8082.3999
$ws.Range("L134").Value = 4509.3333
$ws.Range("M134").Value = -5547.3999
$ws.Range("N134").Value = -9579.3333

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# row 2
$ws.Range("H2").Value = 44997
# row 35
$ws.Range("H35").Value = 5148.923
$ws.Range("I35").Value = 3387.2
$ws.Range("J35").Value = 6250
$ws.Range("K35").Value = 3387.2
$ws.Range("L35").Value = 6250
$ws.Range("M35").Value = -3093.2
$ws.Range("N35").Value = -6838
# row 39
$ws.Range("H39").Value = 10248.625
$ws.Range("I39").Value = 3999.5
$ws.Range("J39").Value = 12331.667
$ws.Range("K39").Value = 3999.5
$ws.Range("L39").Value = 12331.667
$ws.Range("M39").Value = -3608.5
$ws.Range("N39").Value = -13113.667
# row 49
$ws.Range("H49").Value = 10248.625
$ws.Range("I49").Value = 3999.5
$ws.Range("J49").Value = 12331.667
$ws.Range("K49").Value = 3999.5
$ws.Range("L49").Value = 12331.667
$ws.Range("M49").Value = -3817.5
$ws.Range("N49").Value = -12695.667
# row 107
$ws.Range("H107").Value = 1001.0769
$ws.Range("I107").Value = 982
$ws.Range("K107").Value = 982
$ws.Range("M107").Value = 938
# row 141
$ws.Range("H141").Value = 115139.625
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 115139.625
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 115139.625
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -125499.625

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# row 8
$ws.Range("H8").Value = 138
$ws.Range("I8").Value = 138
$ws.Range("K8").Value = 414
$ws.Range("M8").Value = -275
# row 81
$ws.Range("H81").Value = 3334983.2
$ws.Range("J81").Value = 2475
$ws.Range("L81").Value = 7425
$ws.Range("N81").Value = -9671
# row 84
$ws.Range("H84").Value = 3334983.2
$ws.Range("J84").Value = 2475
$ws.Range("L84").Value = 22275
$ws.Range("N84").Value = -33507
# row 113
$ws.Range("H113").Value = 3212.25
$ws.Range("I113").Value = 2966.3333
$ws.Range("J113").Value = 3950
$ws.Range("K113").Value = 8898.999899999999
$ws.Range("L113").Value = 11850
$ws.Range("M113").Value = -6728.999899999999
$ws.Range("N113").Value = -16190
# row 128
$ws.Range("H128").Value = 197015
$ws.Range("I128").Value = 197015
$ws.Range("K128").Value = 591045
$ws.Range("M128").Value = -586065

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 1041927.5
$ws.Range("I2").Value = 3333376.5
$ws.Range("K2").Value = 3333376.5
$ws.Range("M2").Value = -3333263.5
# row 70
$ws.Range("H70").Value = 6876.875
$ws.Range("I70").Value = 5005
$ws.Range("K70").Value = 5005
$ws.Range("M70").Value = -4735
# row 73
$ws.Range("H73").Value = 6876.875
$ws.Range("I73").Value = 5005
$ws.Range("K73").Value = 5005
$ws.Range("M73").Value = -4069
# row 80
$ws.Range("H80").Value = 11637.5
$ws.Range("I80").Value = 15999.333
$ws.Range("J80").Value = 9020.4
$ws.Range("K80").Value = 15999.333
$ws.Range("L80").Value = 9020.4
$ws.Range("M80").Value = -15001.333
$ws.Range("N80").Value = -11016.4
# row 83
$ws.Range("H83").Value = 11637.5
$ws.Range("I83").Value = 15999.333
$ws.Range("J83").Value = 9020.4
$ws.Range("K83").Value = 79996.66500000001
$ws.Range("L83").Value = 45102
$ws.Range("M83").Value = -75004.66500000001
$ws.Range("N83").Value = -55086
# row 102
$ws.Range("H102").Value = 4639.2856
$ws.Range("I102").Value = 3163.6667
$ws.Range("K102").Value = 3163.6667
$ws.Range("M102").Value = -1541.6667

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 1727.625
$ws.Range("I22").Value = 1822.0769
$ws.Range("K22").Value = 1822.0769
$ws.Range("M22").Value = -1527.0769
# row 27
$ws.Range("H27").Value = 1727.625
$ws.Range("I27").Value = 1822.0769
$ws.Range("K27").Value = 1822.0769
$ws.Range("M27").Value = -1715.0769
# row 55
$ws.Range("H55").Value = 908
$ws.Range("J55").Value = 1065.5
$ws.Range("L55").Value = 1065.5
$ws.Range("N55").Value = -1411.5
# row 131
$ws.Range("H131").Value = 69911.5
$ws.Range("J131").Value = 87175
$ws.Range("L131").Value = 87175
$ws.Range("N131").Value = -97255
# row 132
$ws.Range("H132").Value = 5082.6665
$ws.Range("I132").Value = 2749.5
$ws.Range("K132").Value = 8248.5
$ws.Range("M132").Value = -5718.5

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# row 23
$ws.Range("H23").Value = 1667049.9
$ws.Range("I23").Value = 1667049.9
$ws.Range("K23").Value = 1667049.9
$ws.Range("M23").Value = -1666820.9
# row 122
$ws.Range("H122").Value = 8698628
$ws.Range("I122").Value = 1417.6842
$ws.Range("K122").Value = 4253.0526
$ws.Range("M122").Value = -1803.0526
# row 132
$ws.Range("H132").Value = 2270.2727
$ws.Range("I132").Value = 2046.8
$ws.Range("J132").Value = 4505
$ws.Range("K132").Value = 6140.4
$ws.Range("L132").Value = 13515
$ws.Range("M132").Value = -3610.4
$ws.Range("N132").Value = -18575
